$wb = $excel.ActiveWorkbook

# Map of row number -> new "想去人数" (F column) value, shared by both
# the "展览" and "全部类型" sheets, which carry the same event list.
$values = @{
    2 = 13
    3 = 78
    4 = 1464
    5 = 14
    6 = 26
    7 = 9
    8 = 41
    9 = 249
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $values.Keys) {
        $ws.Range("F$row").Value = $values[$row]
    }
}
